# Re-run of the Ptolemy (d,p) spectroscopic-factor output for 116Cd bite 2:
# the usual "shuffle" where levels/strengths shift by a row and the
# re-fit spectroscopic factors / errors change accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.1831305083799375
$ws.Range("E2").Value = 0.008910712656070784
$ws.Range("D3").Value = 0.3488840590562667
$ws.Range("E3").Value = 0.003399601062667641
$ws.Range("D4").Value = 0.03672482619990276
$ws.Range("E4").Value = 0.001010775033024846
$ws.Range("D5").Value = 0.3941170235889925
$ws.Range("D6").Value = 0.05723718536433164
$ws.Range("D7").Value = 0.09335552920724868
$ws.Range("E7").Value = 0.001732370645082966
$ws.Range("D8").Value = 0.0008708892800979852
$ws.Range("E8").Value = 0.0002721529000306204
$ws.Range("D9").Value = 0.02452402671886105
$ws.Range("E9").Value = 0.001042189714480151
$ws.Range("D10").Value = 0.02560799712901147
$ws.Range("E10").Value = 0.001111969223992835
$ws.Range("D11").Value = 0.004987041180272033
$ws.Range("E11").Value = 0.0004482733645188345
$ws.Range("D12").Value = 0.009329231023706615
$ws.Range("E12").Value = 0.0006285652108171123
$ws.Range("D13").Value = 0.5386276983092633
$ws.Range("E13").Value = 0.004397232740232714
$ws.Range("D14").Value = 0.1029674929485241
$ws.Range("E14").Value = 0.007207724506396686
$ws.Range("D15").Value = 0.08568951338016828
$ws.Range("E15").Value = 0.004126522960580352
$ws.Range("D16").Value = 0.211649057127671
$ws.Range("E16").Value = 0.01316225479649695
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 0.1748704491175608
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0.2980911743263773
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 0.009780505139697001
$ws.Range("E19").Value = 0
$ws.Range("B20").Value = 1538.408632604229
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0.08873337631849795
$ws.Range("E20").Value = 0
$ws.Range("B21").Value = 1574
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 0.09007138178732341
$ws.Range("E21").Value = 0.01134757565824547
$ws.Range("B22").Value = 1597
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 0.4693091437181121
$ws.Range("E22").Value = 0.01959991783509517
$ws.Range("B23").Value = 1613.149084845743
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0.003510066718336034
$ws.Range("E23").Value = 0.0003677212752542511
$ws.Range("B24").Value = 1625.344244785309
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 0.00562200978141457
$ws.Range("E24").Value = 0.0008031442544877955
$ws.Range("B25").Value = 1649.346715797556
$ws.Range("D25").Value = 0.02792016584155405
$ws.Range("E25").Value = 0.001025220287939561
$ws.Range("B26").Value = 1732.134406941665
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0.00327370342812974
$ws.Range("E26").Value = 0.0005036466812507292
$ws.Range("B27").Value = 1747.835619913894
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 0.07256872219589559
$ws.Range("E27").Value = 0.007070798572933417
$ws.Range("B28").Value = 1772.725731049579
$ws.Range("D28").Value = 0.007982132275190436
$ws.Range("E28").Value = 0.0005756345390762334
$ws.Range("B29").Value = 1785.161835361124
$ws.Range("D29").Value = 0.04652699544045319
$ws.Range("E29").Value = 0.001464810129856853
$ws.Range("B30").Value = 1795.477475116088
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 0.01233994884654492
$ws.Range("E30").Value = 0.0008897141801584115
$ws.Range("B31").Value = 1805.922029018604
$ws.Range("C31").Value = 4
$ws.Range("D31").Value = 0.06276270414327918
$ws.Range("E31").Value = 0.005671328687645709
$ws.Range("B32").Value = 1818
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0.08290223255493233
$ws.Range("E32").Value = 0.002061105229266274
$ws.Range("D33").Value = 0.01424385565275013
$ws.Range("B34").Value = 1840
$ws.Range("D34").Value = 0.01315969530077994
$ws.Range("D35").Value = 0.002951086682371595
$ws.Range("B36").Value = 1841.298544231364
$ws.Range("C36").Value = 2
$ws.Range("D36").Value = 0.002755942173064324
$ws.Range("E36").Value = 0
$ws.Range("B37").Value = 1851.005074483267
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0.006019806574205516
$ws.Range("E37").Value = 0.0005472551431095923
$ws.Range("B38").Value = 1865.530752678807
$ws.Range("D38").Value = 0.06287310036905422
$ws.Range("E38").Value = 0.002060114189786275
$ws.Range("B39").Value = 1876
$ws.Range("D39").Value = 0.058502030323442
$ws.Range("E39").Value = 0.002107822982421772
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 0.01424527025854633
$ws.Range("E40").Value = 0.00163601143025322
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 0.002012867641817251
$ws.Range("E41").Value = 0
$ws.Range("B42").Value = 1896.361077451411
$ws.Range("D42").Value = 0.2158831510436833
$ws.Range("E42").Value = 0
$ws.Range("B43").Value = 1911.792281651251
$ws.Range("D43").Value = 0.02070507690862477
$ws.Range("E43").Value = 0.001165916952136152
$ws.Range("B44").Value = 1924.67520477692
$ws.Range("D44").Value = 0.02076817857046758
$ws.Range("E44").Value = 0.00581832240106485
$ws.Range("B45").Value = 1933.550109340654
$ws.Range("D45").Value = 0.02145826195326309
$ws.Range("E45").Value = 0.004948786310582415
$ws.Range("B46").Value = 1944.234562301843
$ws.Range("D46").Value = 0.0128703089928519
$ws.Range("E46").Value = 0.002077170122264074
$ws.Range("B47").Value = 1959.440632592807
$ws.Range("D47").Value = 0.01714553141885607
$ws.Range("E47").Value = 0.001391284172413141
$ws.Range("B48").Value = 1970.376107067389
$ws.Range("D48").Value = 0.01989636215762547
$ws.Range("E48").Value = 0.001438786519662999
$ws.Range("B49").Value = 1986.766402625133
$ws.Range("C49").Value = 2
$ws.Range("D49").Value = 0.01285158537403391
$ws.Range("E49").Value = 0.0009504387897195498
$ws.Range("B50").Value = 1995.02814271466
$ws.Range("D50").Value = 0.003744696920085261
$ws.Range("E50").Value = 0.0002130587752902816
$ws.Range("B51").Value = 2014.856993695083
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 0.00202594956888477
$ws.Range("E51").Value = 0.0001199699214421627
$ws.Range("B52").Value = 2048.833064430492
$ws.Range("D52").Value = 0.03897266962965387
$ws.Range("E52").Value = 0.00130752463043211
$ws.Range("B53").Value = 2073.838281975398
$ws.Range("D53").Value = 0.05033354689199342
$ws.Range("E53").Value = 0.00153170556898712
$ws.Range("B54").Value = 2092.576632909581
$ws.Range("D54").Value = 0.03156361021997376
$ws.Range("E54").Value = 0.001241987376362604
$ws.Range("D55").Value = 0.01323984967128855
$ws.Range("B56").Value = 2113.24
$ws.Range("C56").Value = 2
$ws.Range("D56").Value = 0.01161223233179095
$ws.Range("C57").Value = 3
$ws.Range("D57").Value = 0.03558961231351905
$ws.Range("E57").Value = 0.002065468571766731
$ws.Range("D58").Value = 0.01361943528461127
$ws.Range("E58").Value = 0.0009166927595411432
